# "black encoded submissions and added final tuned ones for xgb, cat, lightgbm"
#
# Appends three new submission rows to the bottom of the tracker table:
#   221127_cat_v2data_final      / version_1_1
#   221127_xgb_v2data_final      / version_1_2
#   221127_lightgbm_v2data_final / version_1_3
# all dated 2022-11-27, Hand in = TRUE, By = Maria, then two trailing blank
# (but date-styled) rows, and grows the Tabelle2 table/autofilter and the
# worksheet dimension to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 28
$newRows = 3
$lastNewRow = $lastDataRow + $newRows + 2   # 29..31 data, 32..33 blank -> 33

# --- 1. Copy the date/style formatting of the last existing row (A28) down
#        onto A29:A33 before writing any values, so the new cells inherit
#        style index 1 (the existing "dd/mm/yy" date format) instead of a
#        brand-new style entry. -----------------------------------------
$ws.Range("A" + $lastDataRow).Copy()
for ($r = $lastDataRow + 1; $r -le $lastNewRow; $r++) {
    $ws.Range("A" + $r).PasteSpecial(-4122)   # xlPasteFormats
}

$ramps = @("version_1_1", "version_1_2", "version_1_3")
$names = @("221127_cat_v2data_final", "221127_xgb_v2data_final", "221127_lightgbm_v2data_final")

# --- 2. New submission rows 29-31 ---------------------------------------
# Column C (Name Ramp) is written before column B (Name) so the brand new
# shared strings are appended in the same order as the reference document:
# version_1_1, version_1_2, version_1_3, then the three *_v2data_final
# names. A scratch cell holding a text formula is copied and pasted back
# with "paste values" so the destination cells become plain shared-string
# text cells (no formula, no quote-prefix marker).
for ($i = 0; $i -lt $newRows; $i++) {
    $r = $lastDataRow + 1 + $i
    $ws.Range("Z1").Formula = '="' + $ramps[$i] + '"'
    $ws.Range("Z1").Copy()
    $ws.Range("C" + $r).PasteSpecial(-4163)   # xlPasteValues
}
for ($i = 0; $i -lt $newRows; $i++) {
    $r = $lastDataRow + 1 + $i
    $ws.Range("Z1").Formula = '="' + $names[$i] + '"'
    $ws.Range("Z1").Copy()
    $ws.Range("B" + $r).PasteSpecial(-4163)   # xlPasteValues
}
$ws.Range("Z1").ClearContents()

# Column A (date serial), D (re-uses the existing "TRUE" shared string from
# D4) and E (re-uses the existing "Maria" shared string from E4).
for ($i = 0; $i -lt $newRows; $i++) {
    $r = $lastDataRow + 1 + $i
    $ws.Range("A" + $r).Value = 44892   # 2022-11-27

    $ws.Range("D4").Copy()
    $ws.Range("D" + $r).PasteSpecial(-4163)   # xlPasteValues

    $ws.Range("E4").Copy()
    $ws.Range("E" + $r).PasteSpecial(-4163)   # xlPasteValues
}

# --- 3. Two trailing blank rows (32-33) keep only the copied date style on
#        column A (already pasted in step 1); nothing else to set. -------

# --- 4. Resize the table / autofilter to span the new rows --------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E" + $lastNewRow))

# --- 5. Column B width grows to fit the new, longer submission names ----
$ws.Columns.Item(2).ColumnWidth = 26.33

# --- 6. Move selection similar to the authored workbook -----------------
$ws.Range("B" + ($lastDataRow + $newRows + 1)).Select() | Out-Null
